$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.351.48'

$ws.Range("D3").Value = '1.845.86'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'240.21"
$ws.Range("E5").Value = '  -0.02%  '

$ws.Range("D6").Value = "'0.6275"
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = "'0.07483"
$ws.Range("E8").Value = '  -1.98%  '

$ws.Range("D9").Value = "'0.2897"
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").Value = "'24.36"
$ws.Range("E10").Value = '  -1.41%  '

$ws.Range("D11").Value = "'0.07732"
$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").Value = '1.845.16'
$ws.Range("E12").Value = '  -2.28%  '

$ws.Range("D13").Value = "'4.996"
$ws.Range("E13").Value = '  -0.73%  '

$ws.Range("D14").Value = "'0.6786"
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("D15").Value = "'0.00001033"
$ws.Range("E15").Value = '  -2.90%  '

$ws.Range("D16").Value = "'82.08"
$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("D17").Value = '2.102.98'
$ws.Range("E17").Value = '  -3.84%  '

$ws.Range("D18").Value = "'6.156"

$ws.Range("D19").Value = '29.401.19'
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").Value = "'228.93"
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").Value = "'12.31"
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").Value = "'7.457"
$ws.Range("E23").Value = '  -0.54%  '

$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = "'158.52"
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").Value = "'0.1375"
$ws.Range("E26").Value = '  -0.43%  '

$ws.Range("D27").Value = "'8.398"
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").Value = "'17.50"
$ws.Range("E28").Value = '  -0.96%  '

$ws.Range("D29").Value = "'0.06447"
$ws.Range("E29").Value = '  +14.94%  '

$ws.Range("D30").Value = "'1.385"
$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("D31").Value = "'1.473"
$ws.Range("E31").Value = '  +0.95%  '

$ws.Range("D32").Value = "'4.088"
$ws.Range("E32").Value = '  -0.78%  '

$ws.Range("D33").Value = "'4.056"
$ws.Range("E33").Value = '  -0.50%  '

$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").Value = "'1.139"
$ws.Range("E35").Value = '  -1.97%  '

$ws.Range("D36").Value = "'0.7002"
$ws.Range("E36").Value = '  +1.36%  '

$ws.Range("D37").Value = "'2.577"
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("D38").Value = '1.260.94'
$ws.Range("E38").Value = '  +2.49%  '

$ws.Range("D39").Value = "'2.829"
$ws.Range("E39").Value = '  +4.05%  '

$ws.Range("E40").Value = '  +1.38%  '

$ws.Range("D41").Value = "'6.596"

$ws.Range("D42").Value = "'0.9077"
$ws.Range("E42").Value = '  +0.29%  '

$ws.Range("D43").Value = "'0.9987"
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").Value = '2.007.74'
$ws.Range("E44").Value = '  -18.41%  '

$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").Value = "'66.12"
$ws.Range("E46").Value = '  +0.14%  '

$ws.Range("D47").Value = "'1.745"
$ws.Range("E47").Value = '  +4.05%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = "'7.076"
$ws.Range("E48").Value = '  -1.57%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = "'0.1171"
$ws.Range("E49").Value = '  +2.67%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = "'0.00000000117"
$ws.Range("E50").Value = '  -1.69%  '

$ws.Range("D51").Value = "'8.994"
$ws.Range("E51").Value = '  -0.33%  '
